$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nädal 2")

$ws.Range("B18").Value = 43871
$ws.Range("C18").Value = 0.8125
$ws.Range("G18").Value = "kood"
$ws.Range("H18").Value = "RP with EFCore, juhendi järgi"

$ws.Activate()
$ws.Range("H18").Select()
